$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 3.95
$ws.Range("H2").Value = 1.91
$ws.Range("Q2").Value = 1.98
$ws.Range("I3").Value = 2.32
$ws.Range("P3").Value = 1.89
$ws.Range("F4").Value = 2.3
$ws.Range("G4").Value = 2.32
$ws.Range("H4").Value = 3.45
$ws.Range("I4").Value = 3.55
$ws.Range("J4").Value = 3.5
$ws.Range("AH4").Value = 19
$ws.Range("H5").Value = 3.3
$ws.Range("M5").Value = 1.08
$ws.Range("Q5").Value = 2.02
$ws.Range("T5").Value = 1.81
$ws.Range("AD5").Value = 14.5
$ws.Range("AF5").Value = 16
$ws.Range("AH5").Value = 17.5
$ws.Range("AN5").Value = 22
$ws.Range("K6").Value = 3.75
$ws.Range("O6").Value = 1.38
$ws.Range("AG6").Value = 21
$ws.Range("F7").Value = 2.6
$ws.Range("G7").Value = 2.66
$ws.Range("H7").Value = 2.94
$ws.Range("I7").Value = 3.05
$ws.Range("P7").Value = 1.9
$ws.Range("F8").Value = 5.7
$ws.Range("G8").Value = 6.2
$ws.Range("H8").Value = 1.64
$ws.Range("I8").Value = 1.75
$ws.Range("J8").Value = 3.9
$ws.Range("K8").Value = 4.4
$ws.Range("P8").Value = 2
$ws.Range("Q8").Value = 1.86
$ws.Range("T8").Value = 1.89
$ws.Range("U8").Value = 1.98
$ws.Range("Z8").Value = 11
$ws.Range("AB8").Value = 21
$ws.Range("AE8").Value = 19
$ws.Range("AG8").Value = 24
$ws.Range("AJ8").Value = 190
$ws.Range("AK8").Value = 90
$ws.Range("AN8").Value = 120
$ws.Range("H9").Value = 2.9
$ws.Range("J9").Value = 3.5
$ws.Range("K9").Value = 3.65
$ws.Range("P9").Value = 2.12
$ws.Range("Q9").Value = 1.8
$ws.Range("X9").Value = 20
$ws.Range("Y9").Value = 1000
$ws.Range("AB9").Value = 12.5
$ws.Range("F10").Value = 10.5
$ws.Range("G10").Value = 12.5
$ws.Range("H10").Value = 1.38
$ws.Range("I10").Value = 1.39
$ws.Range("K10").Value = 5.6
$ws.Range("AF10").Value = 120
$ws.Range("H11").Value = 2.02
$ws.Range("I11").Value = 2.08
$ws.Range("J11").Value = 3.85
$ws.Range("I12").Value = 4.9
$ws.Range("P12").Value = 1.86
$ws.Range("U12").Value = 2.02
$ws.Range("AL12").Value = 980
$ws.Range("F13").Value = 3.45
$ws.Range("H13").Value = 2.2
$ws.Range("I13").Value = 2.3
$ws.Range("J13").Value = 3.5
$ws.Range("P13").Value = 1.96
$ws.Range("X13").Value = 17
$ws.Range("F14").Value = 4.4
$ws.Range("G14").Value = 4.7
$ws.Range("H14").Value = 1.98
$ws.Range("I14").Value = 1.99
$ws.Range("K14").Value = 3.65
$ws.Range("M14").Value = 1.09
$ws.Range("P14").Value = 1.76
$ws.Range("Q14").Value = 2.26
$ws.Range("T14").Value = 1.99
$ws.Range("U14").Value = 1.89
$ws.Range("X14").Value = 12.5
$ws.Range("Y14").Value = 7.8
$ws.Range("AA14").Value = 34
$ws.Range("AC14").Value = 1000
$ws.Range("AE14").Value = 29
$ws.Range("AG14").Value = 24
$ws.Range("AH14").Value = 25
$ws.Range("AI14").Value = 60
$ws.Range("F15").Value = 4.7
$ws.Range("G15").Value = 5.4
$ws.Range("H15").Value = 1.81
$ws.Range("I15").Value = 1.96
$ws.Range("J15").Value = 3.5
$ws.Range("P15").Value = 1.81
$ws.Range("Q15").Value = 1.98
$ws.Range("T15").Value = 1.88
$ws.Range("U15").Value = 1.93
$ws.Range("Z15").Value = 12.5
$ws.Range("AA15").Value = 24
$ws.Range("AB15").Value = 17
$ws.Range("AC15").Value = 8.800000000000001
$ws.Range("AE15").Value = 24
$ws.Range("AF15").Value = 44
$ws.Range("AG15").Value = 23
$ws.Range("AI15").Value = 46
$ws.Range("AK15").Value = 85
$ws.Range("AL15").Value = 90
$ws.Range("AM15").Value = 140
$ws.Range("AN15").Value = 110
$ws.Range("AO15").Value = 16
$ws.Range("F16").Value = 3.45
$ws.Range("G16").Value = 3.6
$ws.Range("I16").Value = 2.26
$ws.Range("U16").Value = 2.38
$ws.Range("AH16").Value = 17
$ws.Range("F17").Value = 2.8
$ws.Range("G17").Value = 2.88
$ws.Range("P17").Value = 1.88
$ws.Range("T17").Value = 1.8
$ws.Range("AL17").Value = 160
$ws.Range("F18").Value = 1.68
$ws.Range("T18").Value = 1.83
$ws.Range("F19").Value = 2.56
$ws.Range("J19").Value = 3.3
$ws.Range("AA19").Value = 390
$ws.Range("AB19").Value = 9
$ws.Range("AC19").Value = 7.6
$ws.Range("AE19").Value = 55
$ws.Range("AO19").Value = 180
$ws.Range("N21").Value = 2.84
$ws.Range("P21").Value = 1.61
$ws.Range("U21").Value = 1.89
